$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 60 (shifts existing rows 60-184 down to 61-185,
# pushing the sheet dimension from A1:R184 to A1:R185).
$ws.Rows(60).Insert()

# Populate the newly inserted row with the new "Ajo" price-report record.
$ws.Range("A60").Value2 = 4
$ws.Range("B60").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value2 = "Los Lagos"
$ws.Range("D60").Value2 = 44533
$ws.Range("E60").Value2 = 10
$ws.Range("F60").Value2 = 100112003
$ws.Range("G60").Value2 = "Ajo"
$ws.Range("H60").Value2 = "Chino"
$ws.Range("I60").Value2 = "Primera"
$ws.Range("J60").Value2 = 240
$ws.Range("K60").Value2 = 21000
$ws.Range("L60").Value2 = 21000
$ws.Range("M60").Value2 = 21000
$ws.Range("N60").Value2 = "`$/caja 10 kilos"
$ws.Range("O60").Value2 = "China"
$ws.Range("P60").Value2 = 2100
$ws.Range("Q60").Value2 = 10
$ws.Range("R60").Value2 = "Hortaliza"
